{"js": "// Load all body paragraphs with their text so we can locate the\n// specific ones touched by this edit (letter header content, not the\n// similarly-worded table rows further down in the document).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Update the letter date from September 19 to September 21, 2025.\nconst dateParagraph = paragraphs.items.find(p => p.text === \"September 19, 2025\");\nif (dateParagraph) {\n  dateParagraph.insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2) Split the single-line mailing address (\"929 Story Road, San Jose CA\n// 95122\") into two separate lines: \"929 Story Road\" followed by a new\n// paragraph \"San Jose, CA 95122\". Only the first occurrence (the\n// recipient's mailing address block near the top of the letter) is\n// affected \u2014 the later \"PROPERTY ADDRESS\" table cell keeps its original\n// single-line text.\nconst addressParagraph = paragraphs.items.find(p => p.text === \"929 Story Road, San Jose CA 95122\");\nif (addressParagraph) {\n  addressParagraph.insertText(\"929 Story Road\", Word.InsertLocation.replace);\n  addressParagraph.insertParagraph(\"San Jose, CA 95122\", Word.InsertLocation.after);\n}\n\nawait context.sync();\n\n// 3) Remove the now-unwanted blank \"No Spacing\" paragraph that used to sit\n// directly under \"...Board of Directors\" near the signature block.\nconst boardParagraphs = context.document.body.paragraphs;\nboardParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst boardParagraph = boardParagraphs.items.find(\n  p => p.text === \"Vietnam Town Condominium Owners Association Board of Directors\"\n);\nif (boardParagraph) {\n  const nextParagraph = boardParagraph.getNext();\n  nextParagraph.load(\"text\");\n  await context.sync();\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date from September 19 to September 21, 2025.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"September 19, 2025`r\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2) Split the single-line mailing address (\"929 Story Road, San Jose CA\n# 95122\") into two lines: \"929 Story Road\" followed by a new paragraph\n# \"San Jose, CA 95122\". Only the first matching paragraph (the recipient's\n# address block near the top of the letter) is touched -- the later\n# \"PROPERTY ADDRESS\" table cell keeps its original single-line text.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"929 Story Road, San Jose CA 95122`r\") {\n        $p.Range.Text = \"929 Story Road`rSan Jose, CA 95122\"\n        break\n    }\n}\n\n# 3) Remove the now-unwanted blank \"No Spacing\" paragraph that used to sit\n# directly under \"...Board of Directors\" near the signature block.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Vietnam Town Condominium Owners Association Board of Directors`r\") {\n        $nextPara = $p.Next()\n        if ($nextPara.Range.Text -eq \"`r\") {\n            $nextPara.Range.Delete()\n        }\n        break\n    }\n}\n"}
